# Updated symbol list on Fri Dec 16 20:24:12 UTC 2022 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values to
# the crypto symbol table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    # Force the cell to be treated as text so that numeric-looking
    # strings (e.g. "0.02960") keep their exact original formatting
    # (trailing zeros, leading zeros, etc.) instead of being
    # auto-converted into a floating point number by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    # Restore the cell to the workbook's default (unstyled) appearance
    # so no stray per-cell style is introduced.
    $cell.Style = "Normal"
}

# Row 2 - BNB
Set-TextValue "D2" "243.88"

# Row 3 - OKB
Set-TextValue "D3" "23.59"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.709"

# Row 5 - Cronos
Set-TextValue "D5" "0.05817"

# Row 6 - GateToken
Set-TextValue "D6" "3.408"

# Row 7 - KuCoinToken
Set-TextValue "D7" "6.461"

# Row 8 - FTXToken
Set-TextValue "D8" "1.316"

# Row 9 - MXToken
Set-TextValue "D9" "0.7966"

# Row 11 - MandalaExchangeToken
Set-TextValue "D11" "0.07641"

# Row 12 - LiechtensteinCryptoassetsExchange
Set-TextValue "D12" "0.03208"

# Row 13 - BitrueCoin
Set-TextValue "D13" "0.02960"

# Row 14 - BitMartToken
Set-TextValue "D14" "0.09235"

# Row 15 - BitForexToken
Set-TextValue "D15" "0.001662"

# Row 16 - MCDex
Set-TextValue "D16" "3.254"
Set-TextValue "E16" "15MCDexMCBWorstin24h"

# Row 17 - CoinExToken
Set-TextValue "D17" "0.04754"

# Row 18 - One
Set-TextValue "D18" "0.01247"
Set-TextValue "E18" "17OneONEBestin24h"

# Row 19 - TigerCash
Set-TextValue "D19" "0.006268"

# Row 20 - HotbitToken
Set-TextValue "D20" "0.005399"

# Row 21 - BitKan
Set-TextValue "D21" "0.001062"

# Row 23 - LEO
Set-TextValue "D23" "3.694"

# Row 24 - BTSEToken
Set-TextValue "D24" "2.194"

# Row 25 - BitpandaEcosystemToken
Set-TextValue "D25" "0.3341"

# Row 27 - UpBots
Set-TextValue "D27" "0.0009991"
Set-TextValue "E27" "26UpBotsUBXT"

# Row 40 - IDEX
Set-TextValue "D40" "0.04295"

# Row 41 - KickToken
Set-TextValue "D41" "0.007070"

# Row 43 - BKEXToken
Set-TextValue "D43" "0.1058"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.008871"

# Row 46 - CoinLion
Set-TextValue "D46" "0.00005437"

# Row 48 - CoinbaseStockToken
Set-TextValue "D48" "0.7847"

# Row 49 - BOLO
Set-TextValue "D49" "0.09995"
Set-TextValue "E49" "48BOLOBOLO"
